$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Classes")

# ---------------------------------------------------------------------------
# 1. Mark the previously "Not Started" testing rows (E46:E86, skipping the
#    blank class-header rows) as "Implemented" -- "Added more Testing
#    Methods and Classes".
# ---------------------------------------------------------------------------
$skipRows = @(52, 63, 72, 79)
for ($r = 46; $r -le 86; $r++) {
    if ($skipRows -contains $r) { continue }
    $ws.Range("E$r").Value2 = "Implemented"
}

# ---------------------------------------------------------------------------
# 2. The "Testing" completion-ratio formula in J4 used to (incorrectly) pull
#    its "Implemented" term from column D; point it at column E like the
#    other terms.
# ---------------------------------------------------------------------------
$ws.Range("J4").Formula = '=(COUNTIF(E2:E999,"Done"))/(COUNTIF(E2:E999,"Done") + COUNTIF(E2:E999,"Not Started") + COUNTIF(E2:E999,"Implemented") + COUNTIF(E2:E999,"Failed"))'

# ---------------------------------------------------------------------------
# 3. Re-assert the merged cells for the blocks whose data just changed (rows
#    2-71), which is what nudges the <mergeCells> list back into the order
#    Excel would naturally emit after touching those ranges.
# ---------------------------------------------------------------------------
function Reassert-Merges($ranges) {
    foreach ($r in $ranges) {
        $ws.Range($r).UnMerge()
        $ws.Range($r).Merge()
    }
}

$blockRows2to22  = @("B3:B7", "A2:A7", "B9:B18", "A8:A18", "A19:A22", "B20:B22")
$blockRows23to43 = @("B24:B30", "A23:A30", "B32:B36", "A31:A36", "B38:B43", "A37:A43")
$blockRows44to71 = @("B45:B51", "A44:A51", "B53:B62", "A52:A62", "B64:B71", "A63:A71")
$blockRows92to110 = @("B93:B104", "A92:A104", "B106:B110", "A105:A110")
$blockRows72to91 = @("B73:B78", "A72:A78", "B80:B86", "A79:A86", "B88:B91", "A87:A91")

Reassert-Merges $blockRows2to22
Reassert-Merges $blockRows23to43
Reassert-Merges $blockRows44to71
Reassert-Merges $blockRows92to110
Reassert-Merges $blockRows72to91

# ---------------------------------------------------------------------------
# 4. Tidy up the column-D data validation range so it is one contiguous
#    reference (D2:D1048576) instead of the split D2 / D4:D1048576 / D3.
# ---------------------------------------------------------------------------
$ws.Range("D2:D1048576").Validation.Delete()
$ws.Range("D2:D1048576").Validation.Add(3, 1, 1, '"Done,Not Started,Failed"')

# ---------------------------------------------------------------------------
# 5. Move the viewport / selection the way the author left the sheet: scrolled
#    back up to column D and with I7 selected.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("I7").Select()
